$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) stays text so numeric-looking strings are not
# auto-converted to floating point numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "37.089.44"
$ws.Range("E2").Value = "  -0.80%  "

$ws.Range("D3").Value = "2.024.17"
$ws.Range("E3").Value = "  -1.22%  "

$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "226.17"
$ws.Range("E5").Value = "  -1.25%  "

$ws.Range("D6").Value = "0.607"
$ws.Range("E6").Value = "  -1.43%  "

$ws.Range("E7").Value = "  +0.09%  "

$ws.Range("D8").Value = "55.02"
$ws.Range("E8").Value = "  -3.58%  "

$ws.Range("D9").Value = "0.377"
$ws.Range("E9").Value = "  -2.53%  "

$ws.Range("D10").Value = "0.0783"
$ws.Range("E10").Value = "  -0.73%  "

$ws.Range("E11").Value = "  -4.61%  "

$ws.Range("D12").Value = "2.325.55"
$ws.Range("E12").Value = "  -0.98%  "

$ws.Range("D13").Value = "14.11"
$ws.Range("E13").Value = "  -4.51%  "

$ws.Range("D14").Value = "20.12"
$ws.Range("E14").Value = "  -2.93%  "

$ws.Range("D15").Value = "0.742"
$ws.Range("E15").Value = "  -2.19%  "

$ws.Range("D16").Value = "5.17"
$ws.Range("E16").Value = "  -2.89%  "

$ws.Range("D17").Value = "2.024.96"
$ws.Range("E17").Value = "  -1.18%  "

$ws.Range("D18").Value = "37.010.62"
$ws.Range("E18").Value = "  -0.71%  "

$ws.Range("D19").Value = "6.44"
$ws.Range("E19").Value = "  +5.42%  "

$ws.Range("D20").Value = "68.80"
$ws.Range("E20").Value = "  -1.11%  "

$ws.Range("D21").Value = "0.0₃0815"
$ws.Range("E21").Value = "  -1.70%  "

$ws.Range("D22").Value = "222.79"
$ws.Range("E22").Value = "  -1.50%  "

$ws.Range("E23").Value = "  -0.09%  "

$ws.Range("E24").Value = "  +1.78%  "

$ws.Range("D25").Value = "2.18"
$ws.Range("E25").Value = "  -5.58%  "

$ws.Range("D26").Value = "164.98"
$ws.Range("E26").Value = "  -2.03%  "

$ws.Range("D27").Value = "9.16"
$ws.Range("E27").Value = "  -5.35%  "

$ws.Range("D28").Value = "0.126"
$ws.Range("E28").Value = "  -1.84%  "

$ws.Range("D29").Value = "18.65"
$ws.Range("E29").Value = "  -1.79%  "

$ws.Range("D30").Value = "1.30"
$ws.Range("E30").Value = "  -4.08%  "

$ws.Range("E31").Value = "  -1.54%  "

$ws.Range("D32").Value = "4.49"
$ws.Range("E32").Value = "  -1.40%  "

$ws.Range("D33").Value = "0.0605"
$ws.Range("E33").Value = "  -1.90%  "

$ws.Range("D34").Value = "4.48"
$ws.Range("E34").Value = "  -2.32%  "

$ws.Range("D35").Value = "2.33"
$ws.Range("E35").Value = "  -4.20%  "

$ws.Range("D36").Value = "1.87"
$ws.Range("E36").Value = "  +1.10%  "

$ws.Range("E37").Value = "  +0.31%  "

$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "3.09"
$ws.Range("E38").Value = "  -5.19%  "

$ws.Range("B39").Value = "THORChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D39").Value = "5.47"
$ws.Range("E39").Value = "  +2.89%  "

$ws.Range("D40").Value = "1.462.17"
$ws.Range("E40").Value = "  -0.95%  "

$ws.Range("D41").Value = "0.0213"
$ws.Range("E41").Value = "  -3.98%  "

$ws.Range("D42").Value = "95.26"
$ws.Range("E42").Value = "  -1.15%  "

$ws.Range("E43").Value = "  -3.01%  "

$ws.Range("B44").Value = "Cronos"
$ws.Range("C44").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D44").Value = "0.0911"
$ws.Range("E44").Value = "  -3.60%  "

$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").Value = "16.23"
$ws.Range("E45").Value = "  -5.79%  "

$ws.Range("E46").Value = "  -2.76%  "

$ws.Range("D47").Value = "7.22"
$ws.Range("E47").Value = "  +1.00%  "

$ws.Range("E48").Value = "  -1.55%  "

$ws.Range("D49").Value = "2.93"
$ws.Range("E49").Value = "  +0.24%  "

$ws.Range("D50").Value = "2.214.76"
$ws.Range("E50").Value = "  -1.00%  "

$ws.Range("D51").Value = "3.58"
$ws.Range("E51").Value = "  -9.02%  "

# Restore the default (unstyled) cell style on column D now that the text
# values are set, so the cells match the original formatting.
$ws.Range("D2:D51").Style = "Normal"
